$d = $word.ActiveDocument

# 1. Remove the "IT Support Intern" run entirely (leaving an empty paragraph)
#    Use Find/Replace with wildcard-free plain text; replace with empty string.
$d.Content.Find.Execute("IT Support Intern", $false, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

# 2. Remove the "_GoBack" bookmark if present
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
